$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-5 from 2023-10-08 (45207)
# to 2023-10-09 (45208), keeping existing number formatting/style untouched.
foreach ($row in 2..5) {
    $cell = $ws.Range("C$row")
    $cell.Value2 = $cell.Value2 + 1
}
